$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price and volume(1h) figures for the latest symbol-list refresh.
# Values are written with a leading apostrophe (quote-prefix) so they stay plain text
# (matching the original inlineStr cells) instead of being auto-parsed as number/percent,
# then the cell style is reset to "Normal" so no stray quote-prefix formatting is left behind.

$ws.Range("D2").Value = "'321.65"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-3.21%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'42.94"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-5.75%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.203"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-6.34%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.08177"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-3.65%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'4.316"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-2.56%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.811"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-13.15%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.9492"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-4.20%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1116"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-2.60%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1887"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-2.62%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.09363"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-4.67%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.04618"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-1.95%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'7.455"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-21.24%"
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'-0.31%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001288"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-1.40%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.005762"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-2.11%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.360"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-0.78%"
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'-0.20%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.3365"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'0.31%"
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'0.34%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.2547"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'0.04147"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-0.12%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.001249"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-4.03%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.004286"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-7.02%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.0001200"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-7.89%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0002977"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-0.26%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D38").Value = "'0.02671"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'-2.50%"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05567"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-3.08%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.008152"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'4.06%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.1404"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-2.09%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.006542"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-9.90%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002040"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-3.76%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.007658"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-4.82%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.3202"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-9.83%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006735"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-4.65%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'-0.21%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.003088"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-11.47%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.004097"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'15.85%"
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'-0.21%"
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'-0.21%"
$ws.Range("E51").Style = "Normal"
